$d = $word.ActiveDocument
$f = $d.Fields.Item(1)
$code = $f.Code
$result = $f.Result
$whole = $d.Range($code.Start - 1, $result.End + 1)
$whole.Text = "{m:self.name}"
Write-Output "Set whole done"

$leftover = $d.Range(56, 59)
Write-Output ("Leftover Start/End: " + $leftover.Start + "/" + $leftover.End)
$leftover.Text = ""
Write-Output "Leftover cleared"
